$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Qminus1)
$ws.Range("B2").Value = -0.02498128220729426
$ws.Range("C2").Value = 0.4654300422948336
$ws.Range("D2").Value = 0.403900932749472
$ws.Range("E2").Value = 0.6355320076514416
$ws.Range("F2").Value = 0.6534516332611154

# Row 3 (Q0)
$ws.Range("B3").Value = -0.01304922696893488
$ws.Range("C3").Value = 0.5592693588336196
$ws.Range("D3").Value = 0.5422566057549446
$ws.Range("E3").Value = 0.7363807478166065
$ws.Range("F3").Value = 0.7564405001635656
$ws.Range("G3").Value = 19

# Row 4 (Q1)
$ws.Range("B4").Value = 0.3044081735855184
$ws.Range("C4").Value = 0.5519967918708635
$ws.Range("D4").Value = 0.4951141439019667
$ws.Range("E4").Value = 0.7036434778365864
$ws.Range("F4").Value = 0.6527812235807698
$ws.Range("G4").Value = 18
